$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (B12) with a shared string value, mirroring the pattern
# used by the other rows in the worksheet (B2:B11), and give it a font
# that uses the "Major" font scheme (theme heading font) - this is the
# equivalent of adding a CT_Font with <scheme val="major"/> to the
# style part's font table and referencing it from a new cellXf.
$cell = $ws.Range("B12")
$cell.Value = "FontScheme - Major"
$cell.Font.ThemeFont = 2
